$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New literal values for the H column cells whose value actually changed.
$newValues = @{2="8.1";4="8.1";5="8.1";6="8.1";7="8.41645027969696";8="8.1";9="8.1";10="8.1";11="8.1";13="8.1";14="8.1";15="8.1";16="8.1";17="8.1";19="8.1";20="8.1";21="8.1";22="8.1";23="8.1";25="8.1";26="8.1";27="8.1";28="8.1";29="8.1";31="8.1";32="8.1";33="8.1";34="8.66190818614259";36="8.1";37="8.1";38="8.1";39="8.1";40="8.1";41="8.1";43="8.1";44="8.1";45="8.1";46="8.1";47="8.1";49="8.1";50="8.1";51="8.1";52="8.1";53="8.1";54="8.1";55="8.41645027969696";56="8.1";57="8.1";58="8.1";59="8.1";61="8.1";62="8.1";63="8.1";64="8.1";65="8.1";67="8.1";68="8.1";70="8.1";71="8.1";72="8.1";73="8.41645027969696";74="8.1";75="8.1";76="8.1";77="8.1";78="8.1";79="8.41645027969696";80="8.1";81="8.1";82="8.1";83="8.1";85="8.1";86="8.1";87="8.1";88="8.1";89="8.1";91="8.1";92="8.1";93="8.1";94="8.1";95="8.1";96="8.1";97="8.1";98="8.1";99="8.1";100="8.1";101="8.1";102="8.1";103="8.41645027969696";104="8.1";105="8.1";106="8.1";107="8.1";109="8.1";110="8.1";111="8.1";112="8.1";113="8.1";114="8.1";115="8.1";116="8.1";117="8.1";118="8.1";119="8.1";121="8.1";122="8.1";124="8.1";125="8.1";126="8.1";127="8.41645027969696";128="8.1";130="8.1";131="8.1";132="8.1";133="8.41645027969696";134="8.1";135="8.1";136="8.1";137="8.1";139="8.1";140="8.1";141="8.1";142="8.1";143="8.1";145="8.1";146="8.1";147="8.1";148="8.1";149="8.1";151="8.1";152="8.1";153="8.1";154="8.1";155="8.1";157="8.1";158="8.1";159="8.1";160="8.1";161="8.1";163="8.1";164="8.1";165="8.1";166="8.1";167="8.1";169="8.1";170="8.1";171="8.1";172="8.1";173="8.1";175="8.1";176="8.1";177="8.1";178="8.1";179="8.1";180="8.1";181="8.41645027969696";182="8.1";183="8.1";184="8.1";185="8.1";187="8.1";188="8.1";189="8.1";190="8.1";191="8.1";192="8.1";193="8.41645027969696";194="8.1";195="8.1";196="8.1";197="8.1";199="8.1";200="8.1";201="8.1";202="8.1";203="8.1";205="8.1";206="8.1";207="8.1";208="8.1";209="8.1";211="8.1";212="8.1";215="8.1";216="8.1";217="8.41645027969696";218="8.1";219="8.1";220="8.1";221="8.1";223="8.1";224="8.41645027969696";225="8.1";226="8.1";227="9.36794746093222";229="8.1";230="8.1";231="8.1";232="8.1";233="8.1";235="8.1";236="9.30847351753492";237="8.1";238="8.1";239="8.1";241="8.1";242="8.1";243="8.1";244="8.1";245="8.1";247="8.1";248="8.1";249="8.1";250="8.1";251="8.1";253="8.1";254="8.1";255="8.1";256="8.1";257="8.1";259="8.1";260="8.1";261="8.1";262="8.1";263="8.1";265="8.1";266="9.30847351753492";267="8.1";268="8.1";269="8.1";271="8.1"}

foreach ($key in $newValues.Keys) {
    $ws.Cells.Item([int]$key, 8).Value = [double]$newValues[$key]
}

# Every data row in column H (2..271) gets the new "0.0_ " custom number
# format (Excel appends it once to numFmts / cellXfs and re-uses that xf
# for the whole range).
$ws.Range("H2:H271").NumberFormat = "0.0_ "

# Reproduce the column H/J selection left behind in the saved file.
$r1 = $ws.Range("H1:H1048576")
$r2 = $ws.Range("J1:J1048576")
$u = $excel.Union($r1, $r2)
$u.Select()
$ws.Range("J1").Activate()
